# Fixed a bug in WeightTrigger
# The rows of reel-weight data (rows 3-21) were re-ordered; update each
# row's values in place to match the corrected ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,  501,  9, 52, 30, 75, 45),
    @(4,  801,  3, 67, 65, 52, 45),
    @(5,  901, 16, 15, 45, 60, 60),
    @(6,  401,  9, 48, 67, 75, 45),
    @(7,  201,  9, 30, 15, 45, 30),
    @(8,  701,  3, 90, 45, 97, 15),
    @(9,  1201, 2, 10, 10, 10, 10),
    @(10, 902,  1,  0,  0,  0,  0),
    @(11, 1001,18, 30, 75, 60, 72),
    @(12, 301,  6, 45, 30, 60, 45),
    @(13, 601,  9, 60, 67, 60, 42),
    @(14, 1202, 2, 10, 10, 10, 10),
    @(15, 1203, 3, 15, 15, 15, 15),
    @(16, 2,    0,  2,  2,  2,  2),
    @(17, 3,    0,  3,  3,  3,  3),
    @(18, 802,  0,  4,  5,  4,  0),
    @(19, 1101, 0, 15, 30, 30,  0),
    @(20, 1,    0,  2,  2,  2,  2),
    @(21, 502,  0,  4,  0,  0,  0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
